# This script rearranges the stock-report rows that had become
# out-of-order: for a handful of (B..G) "batch record" rows that share
# the same item, the individual batch rows need to be put back so that
# each row's Batch No / MRP / Qty / Value (columns B,C,D,E,F,G) line up
# the way they originally belonged together.
#
# Implementation: capture the current (B:G) values of every row that
# participates in a re-shuffle, then re-write each target row's B:G
# block using the captured values from its mapped source row. Capturing
# everything up-front avoids clobbering data that a later step still
# needs to read (important for the one 3-row rotation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row -> source row: target row's columns B..G should end up
# holding the values that currently live in the source row's B..G.
$rowMap = @{
    149 = 150
    150 = 149
    161 = 163
    162 = 161
    163 = 162
    183 = 184
    184 = 183
    264 = 265
    265 = 264
    316 = 317
    317 = 318
    318 = 316
    350 = 352
    351 = 351
    352 = 350
    375 = 376
    376 = 375
    382 = 383
    383 = 382
    389 = 390
    390 = 389
    419 = 420
    420 = 419
    421 = 422
    422 = 421
    431 = 432
    432 = 431
    457 = 458
    458 = 457
    583 = 584
    584 = 583
    586 = 587
    587 = 586
    590 = 591
    591 = 590
    593 = 594
    594 = 593
    601 = 602
    602 = 601
    687 = 688
    688 = 687
    709 = 710
    710 = 709
    720 = 721
    721 = 720
    872 = 873
    873 = 872
}

# Snapshot the B:G values (columns 2..7) of every row referenced above
# before any writes happen.
$snapshot = @{}
foreach ($r in $rowMap.Values | Select-Object -Unique) {
    $snapshot[$r] = $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 7)).Value2
}

# Now write each target row's B:G block from the snapshot of its source row.
foreach ($target in $rowMap.Keys) {
    $source = $rowMap[$target]
    $values = $snapshot[$source]
    $destRange = $ws.Range($ws.Cells.Item($target, 2), $ws.Cells.Item($target, 7))
    $destRange.Value2 = $values
}
